# Poland II Liga workbook update
#
# The underlying league database was refreshed and, while re-exporting,
# a handful of match rows that share the same kickoff date/time ended up
# re-ordered relative to each other. The sequential "id" column (A) for
# each row position does NOT move, but everything else describing the
# match (external id, teams, score, result letter and all odds columns)
# moves together with the match it belongs to.
#
# This script re-creates that re-ordering by rotating the row content
# (columns B and E:AD) among the affected row groups, using the values
# captured from the worksheet BEFORE any writes happen (so the cycles
# are applied atomically / consistently).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values travel together with a match record.
# (Column A = sequential row id, C = Div, D = Date are never touched -
#  every row inside a given group already shares the same Div/Date.)
$cols = @(2) + @(5..30)

function Get-RowSnapshot($row) {
    $snap = @{}
    foreach ($c in $cols) {
        $snap[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $snap
}

function Set-RowFromSnapshot($row, $snap) {
    foreach ($c in $cols) {
        $ws.Cells.Item($row, $c).Value = $snap[$c]
    }
}

# Every row that is part of a re-ordering group, and which OTHER row's
# (pre-edit) content it should receive.
$sigma = @{
    2 = 3
    3 = 2

    44 = 46
    45 = 44
    46 = 45

    143 = 144
    144 = 143

    265 = 266
    266 = 265

    271 = 273
    272 = 271
    273 = 272

    307 = 308
    308 = 310
    309 = 311
    310 = 312
    311 = 309
    312 = 307
}

# Snapshot every involved row first, before mutating anything, so that
# source data for later assignments is never already-overwritten data.
$snapshots = @{}
foreach ($row in $sigma.Keys) {
    $snapshots[$row] = Get-RowSnapshot $row
}

# Now write each row's new content from the pre-captured snapshot of its
# source row.
foreach ($row in $sigma.Keys) {
    $srcRow = $sigma[$row]
    Set-RowFromSnapshot $row $snapshots[$srcRow]
}
